# Tournament closure exception handled:
# - Rename sheet to "Validation"
# - Update shared string "Unnamed: 5" header to a single space
# - Recompute/update match results (B:E and H:K) for rows 2-5 and
#   populate the separator column F with the blank header value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "Validation"

# The "Unnamed: 5" column header (F1) becomes a single space
$ws.Range("F1").Value = " "

# Fill column F (separator) for data rows with the same blank value
$ws.Range("F2").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("F5").Value = " "

# Row 2 (team in A2 / G2)
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 2

# Row 3
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6

# Row 4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0

# Row 5
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 4
